# Automatische test-sync: 2025-07-31 21:48:50
# Adds "Testmail #12" (geen geld terug / refund not received) as a new logged
# row to the historical responses sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = $ws.UsedRange.Rows.Count + 1

$ws.Cells.Item($newRow, 1).Value  = "Testmail #12: Ik heb nog geen geld terug."
$ws.Cells.Item($newRow, 2).Value  = "Beste klant,`nBedankt voor uw e-mail. Om uw vraag over het uitblijven van de terugbetaling te kunnen beantwoorden, hebben we wat meer informatie nodig. Kunt u ons alstublieft uw bestelnummer en de datum van de oorspronkelijke aankoop verstrekken? Op die manier kunnen we dit verder voor u onderzoeken en u zo snel mogelijk van dienst zijn.`nMet vriendelijke groet,`n[Bedrijfsnaam] e-mailassistent"
$ws.Cells.Item($newRow, 3).Value  = "Ik heb nog geen geld terug."
$ws.Cells.Item($newRow, 4).Value  = "mailmind.test@zohomail.eu"
$ws.Cells.Item($newRow, 5).Value  = "Retour / Terugbetaling"
$ws.Cells.Item($newRow, 6).Value  = "2025-07-31 21:48:06"
$ws.Cells.Item($newRow, 7).Value  = "Ja"
$ws.Cells.Item($newRow, 8).Value  = "Nee"
$ws.Cells.Item($newRow, 9).Value  = "Ja"
$ws.Cells.Item($newRow, 10).Value = "Nee"

$ws.Rows.Item($newRow).AutoFit()
